# Apply strikethrough formatting to the two paragraphs describing the
# "Micro desafío - Paso 3" and "Micro desafío - code review" text blocks
# (the activity content was marked as completed / struck through).

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Es necesario indicarle al usuario") -or
        $t.StartsWith("Cuando agregamos un formulario")) {
        $p.Range.Font.StrikeThrough = 1
    }
}
